$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename / relabel the experiment run identifiers in column A ---
# (the shared-string table's final order depends on the order these are
# assigned, so the order of the statements below is intentional)
$ws.Cells.Item(4, 1).Value  = "PSO_30_500_04-14-17.48.02_Oscar"
$ws.Cells.Item(5, 1).Value  = "PSO_150_100_04-14-17.48.02_Oscar"
$ws.Cells.Item(12, 1).Value = "PSO_30_100_04-15-01.46.02_Oscar"
$ws.Cells.Item(11, 1).Value = "PSO_30_500_04-15-01.36.02_Aish"
$ws.Cells.Item(13, 1).Value = "PSO_30_500_04-14-18.47.01_Aish"
$ws.Cells.Item(2, 1).Value  = "Oscar PSO_1st_attempt"

# Fix the D12 "30 Error, did 15 instead" note back to a plain particle count
$ws.Cells.Item(12, 4).Value = 30

# Row 14 (the duplicate "Oscar" / 30 / 100 run) is no longer needed
$ws.Rows(14).ClearContents()

# Update the visible selection/scroll position saved with the sheet
$ws.Range("A7").Select() | Out-Null
